$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Fix the "Fixed In"/"Found In" identifier strings: insert a space before the
# trailing "T<n>" suffix (and fix a typo 376 -> 576 on I4).
$ws.Range("H4").Value = "IB*2.0*576 T3"
$ws.Range("I4").Value = "IB*2.0*576 T4"

$ws.Range("H5").Value = "IB*2.0*576 T3"
$ws.Range("I5").Value = "IB*2.0*576 T4"

$ws.Range("H6").Value = "IB*2.0*576 T5"
$ws.Range("I6").Value = "IB*2.0*576 T6"

$ws.Range("H7").Value = "IB*2.0*576 T5"
$ws.Range("I7").Value = "IB*2.0*576 T6"

$ws.Range("H8").Value = "IB*2*576 T6"

$ws.Range("H9").Value = "IB*2*576 T6"

# Remove the stray "z" row that was left below the table.
$ws.Rows("16").Delete()

# Reset the saved selection to the header band.
$ws.Range("A1:J2").Select()
